# "Updating the passwords of users"
#
# - Swap the George / Diego rows (A3 <-> A4) so the order becomes
#   christiano, george, diego.
# - Replace every user's password ("iRequest@1234") with the new
#   password ("iRequest@12").
# - Widen column B to fit the (now shorter) login/password text.
# - Leave the final selection on D17, matching the saved sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-order the user names in column A (George now row 3, Diego now row 4).
$ws.Range("A3").Value = "george.b@mariners.com"
$ws.Range("A4").Value = "diego.maradona@mariners.com"

# New password for every account.
$ws.Range("B2").Value = "iRequest@12"
$ws.Range("B3").Value = "iRequest@12"
$ws.Range("B4").Value = "iRequest@12"

# Widen column B (Password) to show the full values comfortably.
$ws.Columns.Item(2).ColumnWidth = 27.71

# Match the persisted selection/active cell.
$ws.Range("D17").Select()
